$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark changed cells as Text format so values round-trip as strings
$ws.Range("D2:D7").NumberFormat = "@"
$ws.Range("D10:D16").NumberFormat = "@"
$ws.Range("D18:D22").NumberFormat = "@"
$ws.Range("D24:D25").NumberFormat = "@"
$ws.Range("D38:D46").NumberFormat = "@"
$ws.Range("E2:E12").NumberFormat = "@"
$ws.Range("E14:E16").NumberFormat = "@"
$ws.Range("E18:E25").NumberFormat = "@"
$ws.Range("E38:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Step 2: write the new values (as text)
$ws.Range("D2").Value = "291.85"
$ws.Range("E2").Value = "-5.77%"
$ws.Range("G2").Value = "13"
$ws.Range("D3").Value = "40.39"
$ws.Range("E3").Value = "-1.65%"
$ws.Range("G3").Value = "13"
$ws.Range("D4").Value = "5.042"
$ws.Range("E4").Value = "-3.16%"
$ws.Range("G4").Value = "13"
$ws.Range("D5").Value = "0.07369"
$ws.Range("E5").Value = "-4.08%"
$ws.Range("G5").Value = "13"
$ws.Range("D6").Value = "4.288"
$ws.Range("E6").Value = "-0.33%"
$ws.Range("G6").Value = "13"
$ws.Range("D7").Value = "1.566"
$ws.Range("E7").Value = "-8.07%"
$ws.Range("G7").Value = "13"
$ws.Range("E8").Value = "0.55%"
$ws.Range("G8").Value = "13"
$ws.Range("E9").Value = "-4.13%"
$ws.Range("G9").Value = "13"
$ws.Range("D10").Value = "0.1768"
$ws.Range("E10").Value = "-3.39%"
$ws.Range("G10").Value = "13"
$ws.Range("D11").Value = "0.08761"
$ws.Range("E11").Value = "-4.31%"
$ws.Range("G11").Value = "13"
$ws.Range("D12").Value = "0.04197"
$ws.Range("E12").Value = "0.24%"
$ws.Range("G12").Value = "13"
$ws.Range("D13").Value = "0.1053"
$ws.Range("G13").Value = "13"
$ws.Range("D14").Value = "0.001276"
$ws.Range("E14").Value = "1.44%"
$ws.Range("G14").Value = "13"
$ws.Range("D15").Value = "0.005830"
$ws.Range("E15").Value = "0.45%"
$ws.Range("G15").Value = "13"
$ws.Range("D16").Value = "3.413"
$ws.Range("E16").Value = "2.08%"
$ws.Range("G16").Value = "13"
$ws.Range("G17").Value = "13"
$ws.Range("D18").Value = "0.3299"
$ws.Range("E18").Value = "-0.62%"
$ws.Range("G18").Value = "13"
$ws.Range("D19").Value = "7.561"
$ws.Range("E19").Value = "0.99%"
$ws.Range("G19").Value = "13"
$ws.Range("D20").Value = "0.1342"
$ws.Range("E20").Value = "-4.43%"
$ws.Range("G20").Value = "13"
$ws.Range("D21").Value = "0.2863"
$ws.Range("E21").Value = "1.28%"
$ws.Range("G21").Value = "13"
$ws.Range("D22").Value = "0.03835"
$ws.Range("E22").Value = "-4.65%"
$ws.Range("G22").Value = "13"
$ws.Range("E23").Value = "1.06%"
$ws.Range("G23").Value = "13"
$ws.Range("D24").Value = "0.003895"
$ws.Range("E24").Value = "-4.58%"
$ws.Range("G24").Value = "13"
$ws.Range("D25").Value = "0.0001282"
$ws.Range("E25").Value = "-1.57%"
$ws.Range("G25").Value = "13"
$ws.Range("G26").Value = "13"
$ws.Range("G27").Value = "13"
$ws.Range("G28").Value = "13"
$ws.Range("G29").Value = "13"
$ws.Range("G30").Value = "13"
$ws.Range("G31").Value = "13"
$ws.Range("G32").Value = "13"
$ws.Range("G33").Value = "13"
$ws.Range("G34").Value = "13"
$ws.Range("G35").Value = "13"
$ws.Range("G36").Value = "13"
$ws.Range("G37").Value = "13"
$ws.Range("D38").Value = "0.02339"
$ws.Range("E38").Value = "-7.81%"
$ws.Range("G38").Value = "13"
$ws.Range("D39").Value = "0.05050"
$ws.Range("E39").Value = "-5.03%"
$ws.Range("G39").Value = "13"
$ws.Range("D40").Value = "0.007708"
$ws.Range("E40").Value = "-1.86%"
$ws.Range("G40").Value = "13"
$ws.Range("D41").Value = "0.005008"
$ws.Range("E41").Value = "168.63%"
$ws.Range("G41").Value = "13"
$ws.Range("D42").Value = "0.1273"
$ws.Range("E42").Value = "-2.75%"
$ws.Range("G42").Value = "13"
$ws.Range("D43").Value = "0.007370"
$ws.Range("E43").Value = "10.79%"
$ws.Range("G43").Value = "13"
$ws.Range("D44").Value = "0.006960"
$ws.Range("E44").Value = "-5.67%"
$ws.Range("G44").Value = "13"
$ws.Range("D45").Value = "0.3160"
$ws.Range("E45").Value = "2.88%"
$ws.Range("G45").Value = "13"
$ws.Range("D46").Value = "0.00006547"
$ws.Range("E46").Value = "-3.41%"
$ws.Range("G46").Value = "13"
$ws.Range("E47").Value = "-0.03%"
$ws.Range("G47").Value = "13"
$ws.Range("E48").Value = "6.42%"
$ws.Range("G48").Value = "13"
$ws.Range("E49").Value = "35.59%"
$ws.Range("G49").Value = "13"
$ws.Range("E50").Value = "-0.03%"
$ws.Range("G50").Value = "13"
$ws.Range("E51").Value = "-0.03%"
$ws.Range("G51").Value = "13"
